$wb = $excel.ActiveWorkbook

$wsLogin = $wb.Worksheets.Item("SparcsN4Login")
$wsUfv   = $wb.Worksheets.Item("UnitFacilityVisit")
$wsMcr   = $wb.Worksheets.Item("N4MobileCommonRoutines")
$wsCommon = $wb.Worksheets.Item("SparcsN4CommonRoutines")

# ---------------------------------------------------------------------
# The order of writes below mirrors the original authoring order (as
# captured by the workbook's change-tracking logs) so that new shared
# strings are interned in the same sequence: YINV_103, UnitId,
# YINV_104, YINV_105, YINV_106, YINV_107.
# ---------------------------------------------------------------------

# --- N4MobileCommonRoutines: row 4 (YINV_103) -------------------------
$wsMcr.Range("A4").Value = "YINV_103"
$wsMcr.Range("B4").Value = "admin"
$wsMcr.Range("C4").Value = "Admin"
$wsMcr.Range("D4").Value = "OPR1"
$wsMcr.Range("E4").Value = "CPX11"
$wsMcr.Range("F4").Value = "FCY111"
$wsMcr.Range("G4").Value = "YRD1111"
$wsMcr.Range("H4").Value = "Yard Inventory"
$wsMcr.Range("I4").Value = "Query"
$wsMcr.Range("J4").Value = "SBSU1234570"
$wsMcr.Range("K4").Value = "S1"
$wsMcr.Range("L4").Value = "S2"
$wsMcr.Range("M4").Value = "S3"
$wsMcr.Range("N4").Value = "S4"

# --- UnitFacilityVisit: header + row (UnitId) --------------------------
$wsUfv.Range("B1").Value = "UnitId"
$wsUfv.Range("A2").Value = "YINV_103"
$wsUfv.Range("B2").Value = "SBSU1234570"
$wsUfv.Range("G13").Select() | Out-Null

# --- SparcsN4Login: rows 10-12 -----------------------------------------
$wsLogin.Range("A10").Value = "N4MobileCommonRoutines"
$wsLogin.Range("B10").Value = "admin"
$wsLogin.Range("C10").Value = "Admin"
$wsLogin.Range("D10").Value = "OPR1"
$wsLogin.Range("E10").Value = "CPX11"
$wsLogin.Range("F10").Value = "FCY111"
$wsLogin.Range("G10").Value = "YRD1111"
$wsLogin.Range("H10").Value = "Configuration"

$wsLogin.Range("A11").Value = "YINV_101"

$wsLogin.Range("A12").Value = "YINV_103"
$wsLogin.Range("B12").Value = "admin"
$wsLogin.Range("C12").Value = "Admin"
$wsLogin.Range("D12").Value = "OPR1"
$wsLogin.Range("E12").Value = "CPX11"
$wsLogin.Range("F12").Value = "FCY111"
$wsLogin.Range("G12").Value = "YRD1111"
$wsLogin.Range("H12").Value = "Operations"

# --- N4MobileCommonRoutines: row 5 (YINV_104) ---------------------------
$wsMcr.Range("A5").Value = "YINV_104"
$wsMcr.Range("B5").Value = "admin"
$wsMcr.Range("C5").Value = "Admin"
$wsMcr.Range("D5").Value = "OPR1"
$wsMcr.Range("E5").Value = "CPX11"
$wsMcr.Range("F5").Value = "FCY111"
$wsMcr.Range("G5").Value = "YRD1111"
$wsMcr.Range("H5").Value = "Yard Inventory"
$wsMcr.Range("I5").Value = "Query"
$wsMcr.Range("J5").Value = "SBSU1234570"
$wsMcr.Range("K5").Value = "S1"
$wsMcr.Range("L5").Value = "S2"
$wsMcr.Range("M5").Value = "S3"
$wsMcr.Range("N5").Value = "S4"
$wsMcr.Range("O5").Value = "DOOR"
$wsMcr.Range("P5").Value = "DAT"
$wsMcr.Range("Q5").Value = "Major"

# --- N4MobileCommonRoutines: row 6 (YINV_105) ---------------------------
$wsMcr.Range("A6").Value = "YINV_105"
$wsMcr.Range("B6").Value = "admin"
$wsMcr.Range("C6").Value = "Admin"
$wsMcr.Range("D6").Value = "OPR1"
$wsMcr.Range("E6").Value = "CPX11"
$wsMcr.Range("F6").Value = "FCY111"
$wsMcr.Range("G6").Value = "YRD1111"
$wsMcr.Range("H6").Value = "Yard Inventory"
$wsMcr.Range("I6").Value = "Query"
$wsMcr.Range("J6").Value = "SBSU1234570"
$wsMcr.Range("W6").Value = 26000
$wsMcr.Range("X6").Value = 2000

# --- SparcsN4Login: rows 13-14 (YINV_104 / YINV_105) --------------------
$wsLogin.Range("A13").Value = "YINV_104"
$wsLogin.Range("B13").Value = "admin"
$wsLogin.Range("C13").Value = "Admin"
$wsLogin.Range("D13").Value = "OPR1"
$wsLogin.Range("E13").Value = "CPX11"
$wsLogin.Range("F13").Value = "FCY111"
$wsLogin.Range("G13").Value = "YRD1111"
$wsLogin.Range("H13").Value = "Operations"

$wsLogin.Range("A14").Value = "YINV_105"
$wsLogin.Range("B14").Value = "admin"
$wsLogin.Range("C14").Value = "Admin"
$wsLogin.Range("D14").Value = "OPR1"
$wsLogin.Range("E14").Value = "CPX11"
$wsLogin.Range("F14").Value = "FCY111"
$wsLogin.Range("G14").Value = "YRD1111"
$wsLogin.Range("H14").Value = "Operations"

# --- N4MobileCommonRoutines: row 7 (YINV_106) ---------------------------
$wsMcr.Range("A7").Value = "YINV_106"
$wsMcr.Range("B7").Value = "admin"
$wsMcr.Range("C7").Value = "Admin"
$wsMcr.Range("D7").Value = "OPR1"
$wsMcr.Range("E7").Value = "CPX11"
$wsMcr.Range("F7").Value = "FCY111"
$wsMcr.Range("G7").Value = "YRD1111"
$wsMcr.Range("H7").Value = "Yard Inventory"
$wsMcr.Range("I7").Value = "Query"
$wsMcr.Range("J7").Value = "SBSU1234570"
$wsMcr.Range("AB7").Value = 10
$wsMcr.Range("AC7").Value = 11
$wsMcr.Range("AD7").Value = 12
$wsMcr.Range("AE7").Value = 13
$wsMcr.Range("AF7").Value = 10
$wsMcr.Range("AG7").Value = "cm"

# --- SparcsN4Login: row 15 (YINV_106) -----------------------------------
$wsLogin.Range("A15").Value = "YINV_106"
$wsLogin.Range("B15").Value = "admin"
$wsLogin.Range("C15").Value = "Admin"
$wsLogin.Range("D15").Value = "OPR1"
$wsLogin.Range("E15").Value = "CPX11"
$wsLogin.Range("F15").Value = "FCY111"
$wsLogin.Range("G15").Value = "YRD1111"
$wsLogin.Range("H15").Value = "Operations"

# --- SparcsN4Login: row 16 (YINV_107) -----------------------------------
$wsLogin.Range("A16").Value = "YINV_107"
$wsLogin.Range("B16").Value = "admin"
$wsLogin.Range("C16").Value = "Admin"
$wsLogin.Range("D16").Value = "OPR1"
$wsLogin.Range("E16").Value = "CPX11"
$wsLogin.Range("F16").Value = "FCY111"
$wsLogin.Range("G16").Value = "YRD1111"
$wsLogin.Range("H16").Value = "Operations"

$wsLogin.Range("B16:H16").Select() | Out-Null

# --- N4MobileCommonRoutines: row 8 (YINV_107) ---------------------------
$wsMcr.Range("A8").Value = "YINV_107"
$wsMcr.Range("B8").Value = "admin"
$wsMcr.Range("C8").Value = "Admin"
$wsMcr.Range("D8").Value = "OPR1"
$wsMcr.Range("E8").Value = "CPX11"
$wsMcr.Range("F8").Value = "FCY111"
$wsMcr.Range("G8").Value = "YRD1111"
$wsMcr.Range("H8").Value = "Yard Inventory"
$wsMcr.Range("I8").Value = "Query"
$wsMcr.Range("J8").Value = "SBSU1234570"
$wsMcr.Range("AH8").Value = "FOOD"

# ---------------------------------------------------------------------
# Sheet "SparcsN4CommonRoutines" -- cosmetic selection change only
# ---------------------------------------------------------------------
$wsCommon.Range("M10").Select() | Out-Null

# Re-select N4MobileCommonRoutines last so it remains the active tab,
# matching the original workbook's bookView state (activeTab=9).
$wsMcr.Activate() | Out-Null
$wsMcr.Range("E16").Select() | Out-Null
